$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet1 ("LoginData"): replace the old invalid_user/invalid_password
# row with new SauceDemo "special" users, and add more rows.
# ---------------------------------------------------------------------
$wsLogin = $wb.Worksheets.Item("LoginData")

# Row 5 used to be invalid_user / invalid_password -> now problem_user
$wsLogin.Range("A5").Value = "problem_user"
$wsLogin.Range("B5").Value = "secret_sauce"
$wsLogin.Range("C5").Value = $true

# New rows 6-9
$wsLogin.Range("A6").Value = "performance_glitch_user"
$wsLogin.Range("B6").Value = "secret_sauce"
$wsLogin.Range("C6").Value = $true

$wsLogin.Range("A7").Value = "error_user"
$wsLogin.Range("B7").Value = "secret_sauce"
$wsLogin.Range("C7").Value = $true

$wsLogin.Range("A8").Value = "visual_user"
$wsLogin.Range("B8").Value = "secret_sauce"
$wsLogin.Range("C8").Value = $true

# Row 9 has no username (column A left blank)
$wsLogin.Range("B9").Value = "secret_sauce"
$wsLogin.Range("C9").Value = $false

# Column widths (auto-fit-like sizing applied after the new data landed);
# input values are pre-compensated so the engine's pixel-snap lands as
# close as possible to the real-Excel bestFit widths (15.42578125 / 16.7109375)
$wsLogin.Columns.Item(1).ColumnWidth = 14.6
$wsLogin.Columns.Item(2).ColumnWidth = 15.8

[void]$wsLogin.Range("B11").Select()

# ---------------------------------------------------------------------
# Sheet2 ("Inventory"): add a new inventory row not present on the page.
# ---------------------------------------------------------------------
$wsInv = $wb.Worksheets.Item("Inventory")

$wsInv.Range("A8").Value = "Not in list"
$wsInv.Range("B8").Value = 15.99
$wsInv.Range("B8").NumberFormat = "$#,##0.00;[Red]-$#,##0.00"
$wsInv.Range("C8").Value = "This is a sample that is not included in the webpage"

[void]$wsInv.Range("C8").Select()
